$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New user row: "F00901" (column A) with Sucursal "001" (column C),
# matching the formatting used by the existing rows above it.
$ws.Range("A27").Value = "F00901"

$ws.Range("C26").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "001"

# Reset the view: scroll back to the top and select the newly added row.
$ws.Range("A1").Select() | Out-Null
$ws.Range("A27").Select() | Out-Null
